# Auto-generated edit script: apply targeted numeric updates to the
# 'Zalera_Profits' workbook's per-sheet leve-profit tables (ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR), refreshing currentAveragePrice / LevePrice /
# LeveProfit figures pulled by the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 10116.294
$ws.Range("I43").Value = 12299.875
$ws.Range("K43").Value = 12299.875
$ws.Range("M43").Value = -12230.875
# Row 62
$ws.Range("H62").Value = 53336136
$ws.Range("I62").Value = 66667668
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 66667668
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -66667044
$ws.Range("N62").Value = -11248
# Row 65
$ws.Range("H65").Value = 53336136
$ws.Range("I65").Value = 66667668
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 333338340
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -333335220
$ws.Range("N65").Value = -56240
# Row 74
$ws.Range("H74").Value = 9589.223
$ws.Range("I74").Value = 9075.75
$ws.Range("K74").Value = 9075.75
$ws.Range("M74").Value = -8139.75
# Row 76
$ws.Range("H76").Value = 4171008.2
$ws.Range("I76").Value = 5886678
$ws.Range("J76").Value = 4381.7144
$ws.Range("K76").Value = 5886678
$ws.Range("L76").Value = 4381.7144
$ws.Range("M76").Value = -5886363
$ws.Range("N76").Value = -5011.7144
# Row 77
$ws.Range("H77").Value = 9589.223
$ws.Range("I77").Value = 9075.75
$ws.Range("K77").Value = 45378.75
$ws.Range("M77").Value = -40698.75
# Row 79
$ws.Range("H79").Value = 4171008.2
$ws.Range("I79").Value = 5886678
$ws.Range("J79").Value = 4381.7144
$ws.Range("K79").Value = 5886678
$ws.Range("L79").Value = 4381.7144
$ws.Range("M79").Value = -5885586
$ws.Range("N79").Value = -6565.7144
# Row 125
$ws.Range("H125").Value = 1604.875
$ws.Range("I125").Value = 1268
$ws.Range("J125").Value = 2166.3333
$ws.Range("K125").Value = 11412
$ws.Range("L125").Value = 19496.9997
$ws.Range("M125").Value = -8952
$ws.Range("N125").Value = -24416.9997

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 6575440
$ws.Range("I2").Value = 7668013
$ws.Range("K2").Value = 7668013
$ws.Range("M2").Value = -7667900
# Row 32
$ws.Range("H32").Value = 4112.4546
$ws.Range("I32").Value = 4112.4546
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4112.4546
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3825.4546
$ws.Range("N32").ClearContents()
# Row 110
$ws.Range("H110").Value = 17857826
$ws.Range("I110").Value = 17857826
$ws.Range("K110").Value = 17857826
$ws.Range("M110").Value = -17855781
# Row 116
$ws.Range("H116").Value = 6575440
$ws.Range("I116").Value = 7668013
$ws.Range("K116").Value = 7668013
$ws.Range("M116").Value = -7665719

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 6575440
$ws.Range("I3").Value = 7668013
$ws.Range("K3").Value = 7668013
$ws.Range("M3").Value = -7667899

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 28575336
$ws.Range("I31").Value = 76924270
$ws.Range("K31").Value = 76924270
$ws.Range("M31").Value = -76923975
# Row 34
$ws.Range("H34").Value = 28575336
$ws.Range("I34").Value = 76924270
$ws.Range("K34").Value = 76924270
$ws.Range("M34").Value = -76924068
# Row 134
$ws.Range("H134").Value = 5948.393
$ws.Range("I134").Value = 5769
$ws.Range("K134").Value = 17307
$ws.Range("M134").Value = -14772

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 35231876
$ws.Range("I4").Value = 55615890
$ws.Range("J4").Value = 5439862.5
$ws.Range("K4").Value = 166847670
$ws.Range("L4").Value = 16319587.5
$ws.Range("M4").Value = -166847558
$ws.Range("N4").Value = -16319811.5
# Row 87
$ws.Range("H87").Value = 2500
$ws.Range("I87").Value = 1250
$ws.Range("J87").Value = 5000
$ws.Range("K87").Value = 3750
$ws.Range("L87").Value = 15000
$ws.Range("M87").Value = -2502
$ws.Range("N87").Value = -17496
# Row 90
$ws.Range("H90").Value = 2500
$ws.Range("I90").Value = 1250
$ws.Range("J90").Value = 5000
$ws.Range("K90").Value = 11250
$ws.Range("L90").Value = 45000
$ws.Range("M90").Value = -5010
$ws.Range("N90").Value = -57480
# Row 131
$ws.Range("H131").Value = 13337337
$ws.Range("I131").Value = 33334304
$ws.Range("J131").Value = 6026.6665
$ws.Range("K131").Value = 100002912
$ws.Range("L131").Value = 18079.9995
$ws.Range("M131").Value = -99997872
$ws.Range("N131").Value = -28159.9995
# Row 132
$ws.Range("H132").Value = 4490
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 12319.952
$ws.Range("I70").Value = 10761.083
$ws.Range("K70").Value = 10761.083
$ws.Range("M70").Value = -10491.083
# Row 73
$ws.Range("H73").Value = 12319.952
$ws.Range("I73").Value = 10761.083
$ws.Range("K73").Value = 10761.083
$ws.Range("M73").Value = -9825.083000000001
# Row 102
$ws.Range("H102").Value = 1609.0714
$ws.Range("I102").Value = 1322
$ws.Range("K102").Value = 1322
$ws.Range("M102").Value = 300
# Row 113
$ws.Range("H113").Value = 32586.6
$ws.Range("I113").Value = 3233.875
$ws.Range("K113").Value = 3233.875
$ws.Range("M113").Value = -1063.875
# Row 141
$ws.Range("H141").Value = 264714.5
$ws.Range("J141").Value = 264714.5
$ws.Range("L141").Value = 264714.5
$ws.Range("N141").Value = -275074.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 144743.58
$ws.Range("I7").Value = 201759.2
$ws.Range("K7").Value = 201759.2
$ws.Range("M7").Value = -201647.2
# Row 126
$ws.Range("H126").Value = 144743.58
$ws.Range("I126").Value = 201759.2
$ws.Range("K126").Value = 605277.6000000001
$ws.Range("M126").Value = -602807.6000000001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 105494.8
$ws.Range("I126").Value = 116105.336
$ws.Range("J126").Value = 10000
$ws.Range("K126").Value = 348316.008
$ws.Range("L126").Value = 30000
$ws.Range("M126").Value = -345846.008
$ws.Range("N126").Value = -34940
# Row 132
$ws.Range("H132").Value = 7200.773
$ws.Range("I132").Value = 2915.125
$ws.Range("J132").Value = 9649.714
$ws.Range("K132").Value = 8745.375
$ws.Range("L132").Value = 28949.142
$ws.Range("M132").Value = -6215.375
$ws.Range("N132").Value = -34009.142
# Row 140
$ws.Range("H140").Value = 82971.39999999999
$ws.Range("J140").Value = 82971.39999999999
$ws.Range("L140").Value = 82971.39999999999
$ws.Range("N140").Value = -93331.39999999999
